# Update pseudo-F, p-value, and q-value figures in the beta-group
# significance table following refreshed PERMANOVA runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# pseudo-F (column E) updates
$ws.Range("E2").Value  = 4.083406150078668
$ws.Range("E3").Value  = 7.110946498185514
$ws.Range("E4").Value  = 10.26017554980502
$ws.Range("E5").Value  = 5.981288960666151
$ws.Range("E6").Value  = 3.547350313761885
$ws.Range("E7").Value  = 10.91820628338189
$ws.Range("E8").Value  = 5.033447988173603
$ws.Range("E9").Value  = 2.883776193844222
$ws.Range("E10").Value = 2.378368095108534
$ws.Range("E12").Value = 16.32369960621812
$ws.Range("E13").Value = 15.54183971380953
$ws.Range("E14").Value = 9.449050815321602
$ws.Range("E15").Value = 4.109311091116523
$ws.Range("E16").Value = 7.279160948083104
$ws.Range("E17").Value = 6.256488129155477
$ws.Range("E18").Value = 7.656125819316437
$ws.Range("E19").Value = 7.618681469090247
$ws.Range("E20").Value = 6.961038887930381
$ws.Range("E21").Value = 4.308756609755021

# p-value (column F) updates
$ws.Range("F9").Value  = 0.004
$ws.Range("F10").Value = 0.004
$ws.Range("F11").Value = 0.126

# q-value (column G) updates
$ws.Range("G9").Value  = 0.004444444444444445
$ws.Range("G10").Value = 0.004444444444444445
$ws.Range("G11").Value = 0.126
